$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Row 2
Set-TextValue $ws "D2" "304.40"
Set-TextValue $ws "E2" "2.00%"
Set-TextValue $ws "G2" "21"

# Row 3
Set-TextValue $ws "D3" "35.73"
Set-TextValue $ws "E3" "12.38%"
Set-TextValue $ws "G3" "21"

# Row 4
Set-TextValue $ws "D4" "5.050"
Set-TextValue $ws "E4" "1.08%"
Set-TextValue $ws "G4" "21"

# Row 5
Set-TextValue $ws "D5" "0.07798"
Set-TextValue $ws "E5" "0.98%"
Set-TextValue $ws "G5" "21"

# Row 6
Set-TextValue $ws "D6" "2.253"
Set-TextValue $ws "E6" "0.80%"
Set-TextValue $ws "G6" "21"

# Row 7
Set-TextValue $ws "D7" "8.111"
Set-TextValue $ws "E7" "2.62%"
Set-TextValue $ws "G7" "21"

# Row 8
Set-TextValue $ws "D8" "4.040"
Set-TextValue $ws "E8" "6.07%"
Set-TextValue $ws "G8" "21"

# Row 9
Set-TextValue $ws "D9" "0.9299"
Set-TextValue $ws "E9" "0.57%"
Set-TextValue $ws "G9" "21"

# Row 10
Set-TextValue $ws "D10" "0.09568"
Set-TextValue $ws "E10" "-2.30%"
Set-TextValue $ws "G10" "21"

# Row 11
Set-TextValue $ws "D11" "0.1830"
Set-TextValue $ws "E11" "4.62%"
Set-TextValue $ws "G11" "21"

# Row 12
Set-TextValue $ws "D12" "0.08546"
Set-TextValue $ws "E12" "1.35%"
Set-TextValue $ws "G12" "21"

# Row 13
Set-TextValue $ws "E13" "5.01%"
Set-TextValue $ws "G13" "21"

# Row 14
Set-TextValue $ws "D14" "0.09952"
Set-TextValue $ws "E14" "0.79%"
Set-TextValue $ws "G14" "21"

# Row 15
Set-TextValue $ws "D15" "0.001486"
Set-TextValue $ws "E15" "0.88%"
Set-TextValue $ws "G15" "21"

# Row 16
Set-TextValue $ws "D16" "0.005742"
Set-TextValue $ws "E16" "0.22%"
Set-TextValue $ws "G16" "21"

# Row 17
Set-TextValue $ws "E17" "-1.49%"
Set-TextValue $ws "G17" "21"

# Row 18
Set-TextValue $ws "E18" "0.06%"
Set-TextValue $ws "G18" "21"

# Row 19
Set-TextValue $ws "D19" "0.3406"
Set-TextValue $ws "E19" "1.16%"
Set-TextValue $ws "G19" "21"

# Row 20
Set-TextValue $ws "D20" "0.1323"
Set-TextValue $ws "E20" "-0.10%"
Set-TextValue $ws "G20" "21"

# Row 21
Set-TextValue $ws "D21" "4.570"
Set-TextValue $ws "E21" "12.53%"
Set-TextValue $ws "G21" "21"

# Row 22
Set-TextValue $ws "D22" "0.2239"
Set-TextValue $ws "E22" "-1.59%"
Set-TextValue $ws "G22" "21"

# Row 23
Set-TextValue $ws "D23" "0.04682"
Set-TextValue $ws "G23" "21"

# Row 24
Set-TextValue $ws "D24" "0.001244"
Set-TextValue $ws "E24" "2.57%"
Set-TextValue $ws "G24" "21"

# Row 25
Set-TextValue $ws "D25" "0.004538"
Set-TextValue $ws "E25" "3.98%"
Set-TextValue $ws "G25" "21"

# Row 26
Set-TextValue $ws "D26" "0.0001303"
Set-TextValue $ws "E26" "1.18%"
Set-TextValue $ws "G26" "21"

# Row 27
Set-TextValue $ws "E27" "-19.72%"
Set-TextValue $ws "G27" "21"

# Row 28
Set-TextValue $ws "G28" "21"

# Row 29
Set-TextValue $ws "G29" "21"

# Row 30
Set-TextValue $ws "G30" "21"

# Row 31
Set-TextValue $ws "G31" "21"

# Row 32
Set-TextValue $ws "G32" "21"

# Row 33
Set-TextValue $ws "G33" "21"

# Row 34
Set-TextValue $ws "G34" "21"

# Row 35
Set-TextValue $ws "G35" "21"

# Row 36
Set-TextValue $ws "G36" "21"

# Row 37
Set-TextValue $ws "G37" "21"

# Row 38
Set-TextValue $ws "G38" "21"

# Row 39
Set-TextValue $ws "D39" "0.01769"
Set-TextValue $ws "E39" "3.24%"
Set-TextValue $ws "G39" "21"

# Row 40
Set-TextValue $ws "D40" "0.04721"
Set-TextValue $ws "E40" "1.51%"
Set-TextValue $ws "G40" "21"

# Row 41
Set-TextValue $ws "D41" "0.007949"
Set-TextValue $ws "E41" "3.18%"
Set-TextValue $ws "G41" "21"

# Row 42
Set-TextValue $ws "D42" "0.1421"
Set-TextValue $ws "G42" "21"

# Row 43
Set-TextValue $ws "D43" "0.008020"
Set-TextValue $ws "E43" "-17.67%"
Set-TextValue $ws "G43" "21"

# Row 44
Set-TextValue $ws "D44" "0.002228"
Set-TextValue $ws "E44" "12.43%"
Set-TextValue $ws "G44" "21"

# Row 45
Set-TextValue $ws "D45" "0.009114"
Set-TextValue $ws "E45" "-6.17%"
Set-TextValue $ws "G45" "21"

# Row 46
Set-TextValue $ws "D46" "0.00006219"
Set-TextValue $ws "E46" "2.73%"
Set-TextValue $ws "G46" "21"

# Row 47
Set-TextValue $ws "D47" "0.00000000752"
Set-TextValue $ws "E47" "1.17%"
Set-TextValue $ws "G47" "21"

# Row 48
Set-TextValue $ws "D48" "4.049"
Set-TextValue $ws "E48" "52.54%"
Set-TextValue $ws "G48" "21"

# Row 49
Set-TextValue $ws "D49" "0.002697"
Set-TextValue $ws "E49" "36.06%"
Set-TextValue $ws "G49" "21"

# Row 50
Set-TextValue $ws "D50" "0.00002105"
Set-TextValue $ws "E50" "1.17%"
Set-TextValue $ws "G50" "21"

# Row 51
Set-TextValue $ws "D51" "0.0002005"
Set-TextValue $ws "E51" "1.17%"
Set-TextValue $ws "G51" "21"
